# Update slide titles, body text, and citation text for slides 2-11
# (Title 1 = Shapes.Item(1), Content Placeholder 2 = Shapes.Item(2),
#  TextBox 3 (citation) = Shapes.Item(3), citation lives in the 2nd
#  paragraph of the textbox, 1st paragraph is intentionally blank).

$p = $ppt.ActivePresentation

function Set-SlideContent {
    param($SlideIndex, $Title, $Body, $Citation)

    $s = $p.Slides.Item($SlideIndex)

    # Title
    $titleShape = $s.Shapes.Item(1)
    $titleShape.TextFrame.TextRange.Paragraphs(1).Runs(1).Text = $Title

    # Body (bold, 18pt run formatting must be preserved, so only touch the run's text)
    $bodyShape = $s.Shapes.Item(2)
    $bodyShape.TextFrame.TextRange.Paragraphs(1).Runs(1).Text = $Body

    # Citation textbox: paragraph 1 is blank, paragraph 2 holds the citation run
    $citeShape = $s.Shapes.Item(3)
    $originalHeight = $citeShape.Height
    $citeShape.TextFrame.TextRange.Paragraphs(2).Runs(1).Text = $Citation
    # Setting text triggers the shape's auto-fit recalculation; restore the
    # original (unchanged) height so the shape geometry matches the source.
    $citeShape.Height = $originalHeight
}

Set-SlideContent 2 "Slide 1: Introduction to Computing" `
    "Computing is the act of using and developing computer technology. It involves thinking algorithmically and solving problems efficiently. The field of computing encompasses computer science, information technology, and software engineering." `
    "- Smith, J. (2020). Understanding Computing. PublisherX."

Set-SlideContent 3 "Slide 2: History of Computing" `
    "The history of computing dates back to ancient times when the abacus was developed. Over the years, computing devices have evolved from mechanical calculators to modern-day computers and mobile devices." `
    "- Jones, A. (2018). A Brief History of Computing. PublisherY."

Set-SlideContent 4 "Slide 3: Types of Computing Devices" `
    "There are various types of computing devices, including desktop computers, laptops, tablets, and smartphones. Each type of device has its unique features and capabilities." `
    "- Brown, K. (2019). Types of Computing Devices. PublisherZ."

Set-SlideContent 5 "Slide 4: Applications of Computing" `
    "Computing is used in various fields such as healthcare, finance, education, and entertainment. It has revolutionized the way we work, communicate, and access information." `
    "- White, L. (2017). Applications of Computing in the Modern World. PublisherW."

Set-SlideContent 6 "Slide 5: Future Trends in Computing" `
    "The future of computing is predicted to involve artificial intelligence, quantum computing, and edge computing. These technologies have the potential to transform industries and improve efficiency." `
    "- Green, M. (2021). Emerging Trends in Computing. PublisherV."

Set-SlideContent 7 "Slide 6: Importance of Computing Skills" `
    "Having computing skills is essential in today's digital age. It allows individuals to adapt to new technologies, solve complex problems, and innovate in their respective fields." `
    "- Black, S. (2016). The Significance of Computing Skills. PublisherU."

Set-SlideContent 8 "Slide 7: Challenges in Computing" `
    "The field of computing faces challenges such as cybersecurity threats, data privacy issues, and ethical concerns related to artificial intelligence. Addressing these challenges is crucial for the future of computing." `
    "- Red, T. (2015). Challenges in Computing. PublisherT."

Set-SlideContent 9 "Slide 8: Computing Careers" `
    "There are diverse career opportunities in computing, including software development, cybersecurity, data analysis, and network administration. Pursuing a career in computing can lead to rewarding and high-demand job prospects." `
    "- Blue, R. (2014). Careers in Computing. PublisherS."

Set-SlideContent 10 "Slide 9: Computing Education" `
    "Educational programs in computing range from coding bootcamps to computer science degrees. Continuous learning and staying updated with the latest technologies are essential for success in the field of computing." `
    "- Yellow, O. (2013). Computing Education Programs. PublisherR."

Set-SlideContent 11 "Slide 10: Conclusion" `
    "In conclusion, computing plays a vital role in shaping our technological landscape. It is a dynamic and fast-evolving field that offers endless possibilities for innovation and growth." `
    "- Orange, P. (2012). Key Takeaways from Computing. PublisherQ."
